$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = "“트랜스포머를 활용한 자연어 처리”가 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/11/25/%ed%8a%b8%eb%9e%9c%ec%8a%a4%ed%8f%ac%eb%a8%b8%eb%a5%bc-%ed%99%9c%ec%9a%a9%ed%95%9c-%ec%9e%90%ec%97%b0%ec%96%b4-%ec%b2%98%eb%a6%ac%ea%b0%80-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88/"

$ws.Range("D36").Value = "Backbone Network in Deep learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/387"

$ws.Range("D50").Value = "거리 [Kendal tau distance]"
$ws.Range("E50").Value = "http://incredible.egloos.com/7574273"

$ws.Range("D51").Value = "중국 학생들의 낮잠 습관과 학습 성취와의 관계"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EC%A4%91%EA%B5%AD-%ED%95%99%EC%83%9D%EB%93%A4%EC%9D%98-%EB%82%AE%EC%9E%A0-%EC%8A%B5%EA%B4%80%EA%B3%BC-%ED%95%99%EC%8A%B5-%EC%84%B1%EC%B7%A8%EC%99%80%EC%9D%98-%EA%B4%80%EA%B3%84"
